$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "23×32=736" "70×11=770"
Replace-Text "19×67=1273" "33×45=1485"
Replace-Text "99×35=3465" "39×55=2145"
Replace-Text "62×37=2294" "83×21=1743"
Replace-Text "25×32=800" "85×76=6460"
Replace-Text "91×12=1092" "75×18=1350"
Replace-Text "57×88=5016" "65×53=3445"
Replace-Text "86×72=6192" "86×78=6708"
Replace-Text "95×87=8265" "39×31=1209"
Replace-Text "53×20=1060" "57×69=3933"
Replace-Text "93×82=7626" "90×57=5130"
Replace-Text "68×26=1768" "65×35=2275"
Replace-Text "48×65=3120" "86×53=4558"
Replace-Text "98×20=1960" "82×77=6314"
Replace-Text "35×47=1645" "73×41=2993"
Replace-Text "55×16=880" "37×29=1073"
Replace-Text "66×61=4026" "41×93=3813"
Replace-Text "86×82=7052" "66×95=6270"
Replace-Text "74×42=3108" "98×94=9212"
Replace-Text "59×36=2124" "62×23=1426"
Replace-Text "64×86=5504" "19×82=1558"
Replace-Text "20×55=1100" "93×76=7068"
Replace-Text "26×73=1898" "78×71=5538"
Replace-Text "95×86=8170" "55×54=2970"
Replace-Text "60×51=3060" "93×83=7719"

Write-Output "Done applying replacements"
